$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cells below get a new value that looks like a plain decimal number
# (e.g. "13.60", "1.000"). The source workbook stores the Price/Volume
# columns as literal text, but a bare Range.Value assignment lets Excel
# auto-detect such strings as numbers, which silently drops significant
# trailing zeros (e.g. "13.60" -> 13.6). Forcing Text number format first
# keeps the assignment a literal string, matching the source data.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '29.937.65'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '1.882.82'
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '0.7428'
$ws.Range("E5").Value = '  -3.38%  '
$ws.Range("D6").Value = '242.73'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '0.3159'
$ws.Range("E8").Value = '  +0.80%  '
$ws.Range("E9").Value = '  +1.09%  '
$ws.Range("E10").Value = '  -3.26%  '
$ws.Range("D11").Value = '0.08351'
$ws.Range("E11").Value = '  -2.16%  '
$ws.Range("D12").Value = '2.007.69'
$ws.Range("E12").Value = '  +7.07%  '
$ws.Range("D13").Value = '0.7555'
$ws.Range("E13").Value = '  -1.15%  '
$ws.Range("D14").Value = '5.402'
$ws.Range("E14").Value = '  +0.53%  '
$ws.Range("D15").Value = '92.43'
$ws.Range("D16").Value = '6.144'
$ws.Range("E16").Value = '  -0.51%  '
$ws.Range("D17").Value = '30.015.93'
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("D18").Value = '248.88'
$ws.Range("E18").Value = '  +1.76%  '
$ws.Range("D19").Value = '13.60'
$ws.Range("E19").Value = '  -1.34%  '
$ws.Range("D20").Value = '0.000007871'
$ws.Range("E20").Value = '  +0.77%  '
$ws.Range("D21").Value = '2.163.71'
$ws.Range("E21").Value = '  +1.57%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").Value = '8.028'
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").Value = '0.9995'
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").Value = '0.1562'
$ws.Range("E25").Value = '  -3.81%  '
$ws.Range("D26").Value = '9.318'
$ws.Range("E26").Value = '  -1.23%  '
$ws.Range("D27").Value = '166.17'
$ws.Range("E27").Value = '  +1.77%  '
$ws.Range("D28").Value = '18.70'
$ws.Range("E28").Value = '  -0.48%  '
$ws.Range("D29").Value = '2.038'
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").Value = '1.503'
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("D31").Value = '4.602'
$ws.Range("E31").Value = '  +2.35%  '
$ws.Range("D32").Value = '1.537'
$ws.Range("E32").Value = '  -0.26%  '
$ws.Range("D33").Value = '4.221'
$ws.Range("E33").Value = '  +2.58%  '
$ws.Range("D34").Value = '0.05373'
$ws.Range("E34").Value = '  -1.34%  '
$ws.Range("D35").Value = '1.252'
$ws.Range("E35").Value = '  +0.78%  '
$ws.Range("D36").Value = '0.7556'
$ws.Range("E36").Value = '  +1.30%  '
$ws.Range("D37").Value = '1.006'
$ws.Range("E37").Value = '  +0.52%  '
$ws.Range("D38").Value = '2.706'
$ws.Range("E38").Value = '  +0.34%  '
$ws.Range("D39").Value = '0.01968'
$ws.Range("E39").Value = '  +0.92%  '
$ws.Range("D40").Value = '2.760'
$ws.Range("D41").Value = '0.4554'
$ws.Range("E41").Value = '  +1.80%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '6.166'
$ws.Range("E42").Value = '  +1.26%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '1.112.12'
$ws.Range("E43").Value = '  +0.87%  '
$ws.Range("D44").Value = '72.96'
$ws.Range("E44").Value = '  -0.27%  '
$ws.Range("D45").Value = '0.8614'
$ws.Range("E45").Value = '  +0.94%  '
$ws.Range("D46").Value = '104.91'
$ws.Range("E46").Value = '  +1.87%  '
$ws.Range("D47").Value = '1.000'
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").Value = '1.876'
$ws.Range("E48").Value = '  +0.31%  '
$ws.Range("D49").Value = '7.627'
$ws.Range("E49").Value = '  -0.63%  '
$ws.Range("D50").Value = '9.568'
$ws.Range("E50").Value = '  -1.96%  '
$ws.Range("D51").Value = '2.053.83'
$ws.Range("E51").Value = '  +2.45%  '
